$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    # wdFindContinue=1 wrap, MatchCase=$true, MatchWholeWord=$false,
    # Forward=$true, Replace:=wdReplaceAll(2)
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Simple, isolated text swaps. Each of these targets is either the only
#     run in its paragraph, or its neighboring run carries different
#     character formatting, so a plain Find/Replace can't accidentally
#     coalesce it with an adjacent run. ---
Replace-All "9238871" "2130692"
Replace-All "JUKOV DANIL" "АВППВАП"
Replace-All "Студенческий" "Стандартный"
Replace-All "2:45" "02:45"
Replace-All "7:45" "07:45"
Replace-All "т2012" "п235"
Replace-All "Вагон 3" "Вагон 2"
Replace-All "1250" "2500"

# --- The e-registration timestamp run sits immediately after a one-space
#     run that carries the exact same run formatting
#     (<w:b/><w:bCs/><w:lang w:val="en-US"/>). A direct Find/Replace (or a
#     plain Range.Text assignment) on the timestamp would cause the two
#     adjacent, identically-formatted runs to be coalesced into a single
#     run on save; the target revision keeps them as two separate runs
#     (only the timestamp's own <w:t> changes). To avoid the unwanted
#     merge, briefly detune the neighboring space run's bold formatting so
#     it no longer matches while the timestamp text is edited, then
#     restore it afterwards. ---
$oldDate = "13.05.2024 16:23:59"
$newDate = "18.05.2024 11:05:46"

$dateRange = $d.Content.Duplicate
$found = $dateRange.Find.Execute($oldDate)
if ($found) {
    $spaceRange = $d.Range($dateRange.Start - 1, $dateRange.Start)
    $spaceRange.Font.Bold = 0

    $target = $d.Range($dateRange.Start, $dateRange.End)
    $target.Text = $newDate

    $spaceRange = $d.Range($dateRange.Start - 1, $dateRange.Start)
    $spaceRange.Font.Bold = 1
}
